$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.7131737602450967
$ws.Range("C2").Value = -0.7850101286121978
$ws.Range("D2").Value = 0.8011803641307264

$ws.Range("B3").Value = 0.5999845796120017
$ws.Range("C3").Value = -0.8695025230349126
$ws.Range("D3").Value = -0.6211954558013203

$ws.Range("B4").Value = 0.6094797117492496
$ws.Range("C4").Value = -0.7907510129475015
$ws.Range("D4").Value = 0.6473524366391002

$ws.Range("B5").Value = -0.6515268669274742
$ws.Range("C5").Value = -0.5827198377016496
$ws.Range("D5").Value = 0.6887386139170816

$ws.Range("B6").Value = 0.6859373585917701
$ws.Range("C6").Value = -0.671071621034242
$ws.Range("D6").Value = -0.6792891934400925

$ws.Range("B7").Value = -0.7021928898131928
$ws.Range("C7").Value = -0.6944062258610201
$ws.Range("D7").Value = 0.7976176975483849

$ws.Range("B8").Value = -0.7795181946514922
$ws.Range("C8").Value = -0.7061827126860101
$ws.Range("D8").Value = 0.6787788146343691

$ws.Range("B9").Value = 0.7479216463764629
$ws.Range("C9").Value = -0.7850241797102324
$ws.Range("D9").Value = -0.5806547261024267
